$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Anxa1'
$ws.Cells.Item(2, 3).Value = 'Dysf'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 21.39646633333334
$ws.Cells.Item(2, 8).Value = 64.18939900000001
$ws.Cells.Item(2, 9).Value = 0.0721325008796955
$ws.Cells.Item(2, 10).Value = 0.0721325008796955
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 41.30513066666667
$ws.Cells.Item(2, 14).Value = 123.915392
$ws.Cells.Item(2, 15).Value = 0.9052975942275636
$ws.Cells.Item(2, 16).Value = 0.9052975942275634
$ws.Cells.Item(2, 17).Value = 883.7838377032677
$ws.Cells.Item(2, 18).Value = 7954.05453932941
$ws.Cells.Item(2, 19).Value = 0.06530137951200596
$ws.Cells.Item(2, 20).Value = 0.06530137951200593

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Anxa1'
$ws.Cells.Item(3, 3).Value = 'Dysf'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 21.39646633333334
$ws.Cells.Item(3, 8).Value = 64.18939900000001
$ws.Cells.Item(3, 9).Value = 0.0721325008796955
$ws.Cells.Item(3, 10).Value = 0.0721325008796955
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.172896666666667
$ws.Cells.Item(3, 14).Value = 3.51869
$ws.Cells.Item(3, 15).Value = 0.02570674667948099
$ws.Cells.Item(3, 16).Value = 0.02570674667948099
$ws.Cells.Item(3, 17).Value = 25.09584404081223
$ws.Cells.Item(3, 18).Value = 225.8625963673101
$ws.Cells.Item(3, 19).Value = 0.001854291927471772
$ws.Cells.Item(3, 20).Value = 0.001854291927471772

# Row 4
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Anxa1'
$ws.Cells.Item(4, 3).Value = 'Dysf'
$ws.Cells.Item(4, 4).Value = 'M2'
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 21.39646633333334
$ws.Cells.Item(4, 8).Value = 64.18939900000001
$ws.Cells.Item(4, 9).Value = 0.0721325008796955
$ws.Cells.Item(4, 10).Value = 0.0721325008796955
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.05821433333333333
$ws.Cells.Item(4, 14).Value = 0.174643
$ws.Cells.Item(4, 15).Value = 0.001275901929509164
$ws.Cells.Item(4, 16).Value = 0.001275901929509163
$ws.Cells.Item(4, 17).Value = 1.245581023284111
$ws.Cells.Item(4, 18).Value = 11.210229209557
$ws.Cells.Item(4, 19).Value = 0.00009203399705272492
$ws.Cells.Item(4, 20).Value = 0.00009203399705272491

# Row 5
$ws.Cells.Item(5, 1).Value = 'ECs'
$ws.Cells.Item(5, 2).Value = 'Anxa1'
$ws.Cells.Item(5, 3).Value = 'Dysf'
$ws.Cells.Item(5, 4).Value = 'sCs'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 21.39646633333334
$ws.Cells.Item(5, 8).Value = 64.18939900000001
$ws.Cells.Item(5, 9).Value = 0.0721325008796955
$ws.Cells.Item(5, 10).Value = 0.0721325008796955
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.089783333333333
$ws.Cells.Item(5, 14).Value = 9.269349999999999
$ws.Cells.Item(5, 15).Value = 0.06771975716344637
$ws.Cells.Item(5, 16).Value = 0.06771975716344636
$ws.Cells.Item(5, 17).Value = 66.11044506896111
$ws.Cells.Item(5, 18).Value = 594.99400562065
$ws.Cells.Item(5, 19).Value = 0.004884795443165061
$ws.Cells.Item(5, 20).Value = 0.00488479544316506

# Row 6
$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Anxa1'
$ws.Cells.Item(6, 3).Value = 'Dysf'
$ws.Cells.Item(6, 4).Value = 'ECs'
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 130.868154
$ws.Cells.Item(6, 8).Value = 392.604462
$ws.Cells.Item(6, 9).Value = 0.441187207572817
$ws.Cells.Item(6, 10).Value = 0.441187207572817
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 41.30513066666667
$ws.Cells.Item(6, 14).Value = 123.915392
$ws.Cells.Item(6, 15).Value = 0.9052975942275636
$ws.Cells.Item(6, 16).Value = 0.9052975942275634
$ws.Cells.Item(6, 17).Value = 5405.526201075457
$ws.Cells.Item(6, 18).Value = 48649.73580967911
$ws.Cells.Item(6, 19).Value = 0.399405717619648
$ws.Cells.Item(6, 20).Value = 0.3994057176196479

# Row 7
$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Anxa1'
$ws.Cells.Item(7, 3).Value = 'Dysf'
$ws.Cells.Item(7, 4).Value = 'FAPs'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 130.868154
$ws.Cells.Item(7, 8).Value = 392.604462
$ws.Cells.Item(7, 9).Value = 0.441187207572817
$ws.Cells.Item(7, 10).Value = 0.441187207572817
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.172896666666667
$ws.Cells.Item(7, 14).Value = 3.51869
$ws.Cells.Item(7, 15).Value = 0.02570674667948099
$ws.Cells.Item(7, 16).Value = 0.02570674667948099
$ws.Cells.Item(7, 17).Value = 153.49482159942
$ws.Cells.Item(7, 18).Value = 1381.45339439478
$ws.Cells.Item(7, 19).Value = 0.01134148778330201
$ws.Cells.Item(7, 20).Value = 0.011341487783302

# Row 8
$ws.Cells.Item(8, 1).Value = 'FAPs'
$ws.Cells.Item(8, 2).Value = 'Anxa1'
$ws.Cells.Item(8, 3).Value = 'Dysf'
$ws.Cells.Item(8, 4).Value = 'M2'
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 130.868154
$ws.Cells.Item(8, 8).Value = 392.604462
$ws.Cells.Item(8, 9).Value = 0.441187207572817
$ws.Cells.Item(8, 10).Value = 0.441187207572817
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.05821433333333333
$ws.Cells.Item(8, 14).Value = 0.174643
$ws.Cells.Item(8, 15).Value = 0.001275901929509164
$ws.Cells.Item(8, 16).Value = 0.001275901929509163
$ws.Cells.Item(8, 17).Value = 7.618402339674001
$ws.Cells.Item(8, 18).Value = 68.565621057066
$ws.Cells.Item(8, 19).Value = 0.0005629116094169171
$ws.Cells.Item(8, 20).Value = 0.000562911609416917

# Row 9
$ws.Cells.Item(9, 1).Value = 'FAPs'
$ws.Cells.Item(9, 2).Value = 'Anxa1'
$ws.Cells.Item(9, 3).Value = 'Dysf'
$ws.Cells.Item(9, 4).Value = 'sCs'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 130.868154
$ws.Cells.Item(9, 8).Value = 392.604462
$ws.Cells.Item(9, 9).Value = 0.441187207572817
$ws.Cells.Item(9, 10).Value = 0.441187207572817
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.089783333333333
$ws.Cells.Item(9, 14).Value = 9.269349999999999
$ws.Cells.Item(9, 15).Value = 0.06771975716344637
$ws.Cells.Item(9, 16).Value = 0.06771975716344636
$ws.Cells.Item(9, 17).Value = 404.3542410933
$ws.Cells.Item(9, 18).Value = 3639.1881698397
$ws.Cells.Item(9, 19).Value = 0.02987709056045018
$ws.Cells.Item(9, 20).Value = 0.02987709056045017

# Row 10
$ws.Cells.Item(10, 1).Value = 'M2'
$ws.Cells.Item(10, 2).Value = 'Anxa1'
$ws.Cells.Item(10, 3).Value = 'Dysf'
$ws.Cells.Item(10, 4).Value = 'ECs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 75.02619166666666
$ws.Cells.Item(10, 8).Value = 225.078575
$ws.Cells.Item(10, 9).Value = 0.2529308696158396
$ws.Cells.Item(10, 10).Value = 0.2529308696158397
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 41.30513066666667
$ws.Cells.Item(10, 14).Value = 123.915392
$ws.Cells.Item(10, 15).Value = 0.9052975942275636
$ws.Cells.Item(10, 16).Value = 0.9052975942275634
$ws.Cells.Item(10, 17).Value = 3098.966650214044
$ws.Cells.Item(10, 18).Value = 27890.6998519264
$ws.Cells.Item(10, 19).Value = 0.2289777077691052
$ws.Cells.Item(10, 20).Value = 0.2289777077691051

# Row 11
$ws.Cells.Item(11, 1).Value = 'M2'
$ws.Cells.Item(11, 2).Value = 'Anxa1'
$ws.Cells.Item(11, 3).Value = 'Dysf'
$ws.Cells.Item(11, 4).Value = 'FAPs'
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 75.02619166666666
$ws.Cells.Item(11, 8).Value = 225.078575
$ws.Cells.Item(11, 9).Value = 0.2529308696158396
$ws.Cells.Item(11, 10).Value = 0.2529308696158397
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.172896666666667
$ws.Cells.Item(11, 14).Value = 3.51869
$ws.Cells.Item(11, 15).Value = 0.02570674667948099
$ws.Cells.Item(11, 16).Value = 0.02570674667948099
$ws.Cells.Item(11, 17).Value = 87.99797011852777
$ws.Cells.Item(11, 18).Value = 791.9817310667501
$ws.Cells.Item(11, 19).Value = 0.006502029792635224
$ws.Cells.Item(11, 20).Value = 0.006502029792635225

# Row 12
$ws.Cells.Item(12, 1).Value = 'M2'
$ws.Cells.Item(12, 2).Value = 'Anxa1'
$ws.Cells.Item(12, 3).Value = 'Dysf'
$ws.Cells.Item(12, 4).Value = 'M2'
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 75.02619166666666
$ws.Cells.Item(12, 8).Value = 225.078575
$ws.Cells.Item(12, 9).Value = 0.2529308696158396
$ws.Cells.Item(12, 10).Value = 0.2529308696158397
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.05821433333333333
$ws.Cells.Item(12, 14).Value = 0.174643
$ws.Cells.Item(12, 15).Value = 0.001275901929509164
$ws.Cells.Item(12, 16).Value = 0.001275901929509163
$ws.Cells.Item(12, 17).Value = 4.367599730413889
$ws.Cells.Item(12, 18).Value = 39.308397573725
$ws.Cells.Item(12, 19).Value = 0.0003227149845752804
$ws.Cells.Item(12, 20).Value = 0.0003227149845752804

# Row 13
$ws.Cells.Item(13, 1).Value = 'M2'
$ws.Cells.Item(13, 2).Value = 'Anxa1'
$ws.Cells.Item(13, 3).Value = 'Dysf'
$ws.Cells.Item(13, 4).Value = 'sCs'
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 75.02619166666666
$ws.Cells.Item(13, 8).Value = 225.078575
$ws.Cells.Item(13, 9).Value = 0.2529308696158396
$ws.Cells.Item(13, 10).Value = 0.2529308696158397
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.089783333333333
$ws.Cells.Item(13, 14).Value = 9.269349999999999
$ws.Cells.Item(13, 15).Value = 0.06771975716344637
$ws.Cells.Item(13, 16).Value = 0.06771975716344636
$ws.Cells.Item(13, 17).Value = 231.8146765751389
$ws.Cells.Item(13, 18).Value = 2086.33208917625
$ws.Cells.Item(13, 19).Value = 0.01712841706952397
$ws.Cells.Item(13, 20).Value = 0.01712841706952397

# Row 14
$ws.Cells.Item(14, 1).Value = 'sCs'
$ws.Cells.Item(14, 2).Value = 'Anxa1'
$ws.Cells.Item(14, 3).Value = 'Dysf'
$ws.Cells.Item(14, 4).Value = 'ECs'
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 69.33645133333333
$ws.Cells.Item(14, 8).Value = 208.009354
$ws.Cells.Item(14, 9).Value = 0.2337494219316478
$ws.Cells.Item(14, 10).Value = 0.2337494219316478
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 41.30513066666667
$ws.Cells.Item(14, 14).Value = 123.915392
$ws.Cells.Item(14, 15).Value = 0.9052975942275636
$ws.Cells.Item(14, 16).Value = 0.9052975942275634
$ws.Cells.Item(14, 17).Value = 2863.951182286308
$ws.Cells.Item(14, 18).Value = 25775.56064057677
$ws.Cells.Item(14, 19).Value = 0.2116127893268044
$ws.Cells.Item(14, 20).Value = 0.2116127893268044

# Row 15
$ws.Cells.Item(15, 1).Value = 'sCs'
$ws.Cells.Item(15, 2).Value = 'Anxa1'
$ws.Cells.Item(15, 3).Value = 'Dysf'
$ws.Cells.Item(15, 4).Value = 'FAPs'
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 69.33645133333333
$ws.Cells.Item(15, 8).Value = 208.009354
$ws.Cells.Item(15, 9).Value = 0.2337494219316478
$ws.Cells.Item(15, 10).Value = 0.2337494219316478
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.172896666666667
$ws.Cells.Item(15, 14).Value = 3.51869
$ws.Cells.Item(15, 15).Value = 0.02570674667948099
$ws.Cells.Item(15, 16).Value = 0.02570674667948099
$ws.Cells.Item(15, 17).Value = 81.32449264736222
$ws.Cells.Item(15, 18).Value = 731.92043382626
$ws.Cells.Item(15, 19).Value = 0.006008937176071988
$ws.Cells.Item(15, 20).Value = 0.006008937176071988

# Row 16
$ws.Cells.Item(16, 1).Value = 'sCs'
$ws.Cells.Item(16, 2).Value = 'Anxa1'
$ws.Cells.Item(16, 3).Value = 'Dysf'
$ws.Cells.Item(16, 4).Value = 'M2'
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 69.33645133333333
$ws.Cells.Item(16, 8).Value = 208.009354
$ws.Cells.Item(16, 9).Value = 0.2337494219316478
$ws.Cells.Item(16, 10).Value = 0.2337494219316478
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.05821433333333333
$ws.Cells.Item(16, 14).Value = 0.174643
$ws.Cells.Item(16, 15).Value = 0.001275901929509164
$ws.Cells.Item(16, 16).Value = 0.001275901929509163
$ws.Cells.Item(16, 17).Value = 4.036375290069111
$ws.Cells.Item(16, 18).Value = 36.327377610622
$ws.Cells.Item(16, 19).Value = 0.000298241338464241
$ws.Cells.Item(16, 20).Value = 0.0002982413384642409

# Row 17
$ws.Cells.Item(17, 1).Value = 'sCs'
$ws.Cells.Item(17, 2).Value = 'Anxa1'
$ws.Cells.Item(17, 3).Value = 'Dysf'
$ws.Cells.Item(17, 4).Value = 'sCs'
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 69.33645133333333
$ws.Cells.Item(17, 8).Value = 208.009354
$ws.Cells.Item(17, 9).Value = 0.2337494219316478
$ws.Cells.Item(17, 10).Value = 0.2337494219316478
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 3.089783333333333
$ws.Cells.Item(17, 14).Value = 9.269349999999999
$ws.Cells.Item(17, 15).Value = 0.06771975716344637
$ws.Cells.Item(17, 16).Value = 0.06771975716344636
$ws.Cells.Item(17, 17).Value = 214.2346117222111
$ws.Cells.Item(17, 18).Value = 1928.1115054999
$ws.Cells.Item(17, 19).Value = 0.01582945409030715
$ws.Cells.Item(17, 20).Value = 0.01582945409030715
